$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.Value = "'51.720.63"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.Value = "'  -1.13%  "
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.Value = "'2.781.52"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.Value = "'  -1.66%  "
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.Value = "'  -0.02%  "
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.Value = "'357.11"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.Value = "'  +0.20%  "
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.Value = "'109.49"
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.Value = "'  -2.61%  "
$r.Style = "Normal"

$r = $ws.Range("D7")
$r.Value = "'0.554"
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.Value = "'  -3.36%  "
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.Value = "'  -0.01%  "
$r.Style = "Normal"

$r = $ws.Range("D9")
$r.Value = "'0.587"
$r.Style = "Normal"
$r = $ws.Range("E9")
$r.Value = "'  -2.19%  "
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.Value = "'39.69"
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.Value = "'  +3.68%  "
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.Value = "'0.0846"
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.Value = "'19.56"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.Value = "'  -1.92%  "
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.Value = "'7.61"
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.Value = "'  -2.30%  "
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.Value = "'3.220.08"
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.Value = "'  -1.67%  "
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.Value = "'2.767.39"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.Value = "'  -2.51%  "
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.Value = "'  +0.54%  "
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.Value = "'51.713.71"
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.Value = "'  -1.02%  "
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.Value = "'7.57"
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.Value = "'  +0.95%  "
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.Value = "'  -3.51%  "
$r.Style = "Normal"

$r = $ws.Range("D21")
$r.Value = "'13.22"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.Value = "'  -2.38%  "
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.Value = "'  -2.76%  "
$r.Style = "Normal"

$r = $ws.Range("D23")
$r.Value = "'70.22"
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.Value = "'267.66"
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.Value = "'  -1.29%  "
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.Value = "'  -2.56%  "
$r.Style = "Normal"

$r = $ws.Range("D26")
$r.Value = "'26.35"
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.Value = "'  -2.39%  "
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.Value = "'  +17.10%  "
$r.Style = "Normal"

$r = $ws.Range("D28")
$r.Value = "'0.999"
$r.Style = "Normal"
$r = $ws.Range("E28")
$r.Value = "'  -0.09%  "
$r.Style = "Normal"

$r = $ws.Range("B29")
$r.Value = "'Toncoin"
$r.Style = "Normal"
$r = $ws.Range("C29")
$r.Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$r.Style = "Normal"
$r = $ws.Range("D29")
$r.Value = "'2.28"
$r.Style = "Normal"
$r = $ws.Range("E29")
$r.Value = "'  +0.96%  "
$r.Style = "Normal"

$r = $ws.Range("B30")
$r.Value = "'Cosmos"
$r.Style = "Normal"
$r = $ws.Range("C30")
$r.Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$r.Style = "Normal"
$r = $ws.Range("D30")
$r.Value = "'10.19"
$r.Style = "Normal"
$r = $ws.Range("E30")
$r.Value = "'  -1.56%  "
$r.Style = "Normal"

$r = $ws.Range("D31")
$r.Value = "'6.17"
$r.Style = "Normal"
$r = $ws.Range("E31")
$r.Value = "'  +3.87%  "
$r.Style = "Normal"

$r = $ws.Range("D32")
$r.Value = "'35.14"
$r.Style = "Normal"
$r = $ws.Range("E32")
$r.Value = "'  +0.10%  "
$r.Style = "Normal"

$r = $ws.Range("D33")
$r.Value = "'51.98"
$r.Style = "Normal"
$r = $ws.Range("E33")
$r.Value = "'  -1.20%  "
$r.Style = "Normal"

$r = $ws.Range("D34")
$r.Value = "'0.0448"
$r.Style = "Normal"
$r = $ws.Range("E34")
$r.Value = "'  -8.21%  "
$r.Style = "Normal"

$r = $ws.Range("D35")
$r.Value = "'0.0835"
$r.Style = "Normal"
$r = $ws.Range("E35")
$r.Value = "'  -2.52%  "
$r.Style = "Normal"

$r = $ws.Range("D36")
$r.Value = "'5.22"
$r.Style = "Normal"
$r = $ws.Range("E36")
$r.Value = "'  -6.65%  "
$r.Style = "Normal"

$r = $ws.Range("D37")
$r.Value = "'1.00"
$r.Style = "Normal"

$r = $ws.Range("D38")
$r.Value = "'18.86"
$r.Style = "Normal"
$r = $ws.Range("E38")
$r.Value = "'  +1.99%  "
$r.Style = "Normal"

$r = $ws.Range("D39")
$r.Value = "'3.14"
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.Value = "'  -4.23%  "
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.Value = "'  -4.31%  "
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.Value = "'2.54"
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.Value = "'  -0.08%  "
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.Value = "'  -3.04%  "
$r.Style = "Normal"

$r = $ws.Range("B43")
$r.Value = "'Monero"
$r.Style = "Normal"
$r = $ws.Range("C43")
$r.Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$r.Style = "Normal"
$r = $ws.Range("D43")
$r.Value = "'119.84"
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.Value = "'  -6.03%  "
$r.Style = "Normal"

$r = $ws.Range("B44")
$r.Value = "'WEMIXToken"
$r.Style = "Normal"
$r = $ws.Range("C44")
$r.Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$r.Style = "Normal"
$r = $ws.Range("D44")
$r.Value = "'2.20"
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.Value = "'  -3.47%  "
$r.Style = "Normal"

$r = $ws.Range("D45")
$r.Value = "'21.93"
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.Value = "'  -6.05%  "
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.Value = "'2.084.36"
$r.Style = "Normal"
$r = $ws.Range("E46")
$r.Value = "'  -0.37%  "
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.Value = "'  -3.14%  "
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.Value = "'  -0.31%  "
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.Value = "'  -2.25%  "
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.Value = "'  -6.62%  "
$r.Style = "Normal"

$r = $ws.Range("E51")
$r.Value = "'  -1.95%  "
$r.Style = "Normal"
